# Update the "想去人数" (F column) counts on each sheet to reflect the
# latest generated output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 274
$ws.Range("F4").Value  = 1114
$ws.Range("F5").Value  = 2711
$ws.Range("F7").Value  = 689
$ws.Range("F9").Value  = 260
$ws.Range("F10").Value = 192
$ws.Range("F12").Value = 99
$ws.Range("F13").Value = 127
$ws.Range("F14").Value = 1593
$ws.Range("F15").Value = 305
$ws.Range("F18").Value = 255

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 27
$ws.Range("F6").Value  = 15
$ws.Range("F10").Value = 22
$ws.Range("F12").Value = 47
$ws.Range("F17").Value = 6
$ws.Range("F19").Value = 51

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 253

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 253
$ws.Range("F9").Value  = 27
$ws.Range("F11").Value = 274
$ws.Range("F12").Value = 1114
$ws.Range("F13").Value = 15
$ws.Range("F16").Value = 2711
$ws.Range("F19").Value = 22
$ws.Range("F21").Value = 47
$ws.Range("F22").Value = 689
$ws.Range("F24").Value = 260
$ws.Range("F26").Value = 192
$ws.Range("F28").Value = 99
$ws.Range("F29").Value = 127
$ws.Range("F31").Value = 1593
$ws.Range("F32").Value = 305
$ws.Range("F37").Value = 6
$ws.Range("F39").Value = 51
$ws.Range("F43").Value = 255

$wb.Save()
